$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = 10.93086033333333
$ws.Cells.Item(2, 8).Value = 32.792581
$ws.Cells.Item(2, 9).Value = 0.02611891973042464
$ws.Cells.Item(2, 10).Value = 0.02622656909968252
$ws.Cells.Item(2, 11).Value = 1
$ws.Cells.Item(2, 12).Value = 0.3333333333333333
$ws.Cells.Item(2, 13).Value = 0.1375686666666667
$ws.Cells.Item(2, 14).Value = 0.412706
$ws.Cells.Item(2, 15).Value = 0.2896572731203081
$ws.Cells.Item(2, 16).Value = 0.2896572731203081
$ws.Cells.Item(2, 17).Value = 1.503743881576222
$ws.Cells.Item(2, 18).Value = 13.533694934186
$ws.Cells.Item(2, 19).Value = 0.007565535065963015
$ws.Cells.Item(2, 20).Value = 0.007596716488715374

$ws.Cells.Item(3, 7).Value = 10.93086033333333
$ws.Cells.Item(3, 8).Value = 32.792581
$ws.Cells.Item(3, 9).Value = 0.02611891973042464
$ws.Cells.Item(3, 10).Value = 0.02622656909968252
$ws.Cells.Item(3, 11).Value = 2
$ws.Cells.Item(3, 12).Value = 0.6666666666666666
$ws.Cells.Item(3, 13).Value = 0.3373673333333334
$ws.Cells.Item(3, 14).Value = 1.012102
$ws.Cells.Item(3, 15).Value = 0.7103427268796919
$ws.Cells.Item(3, 16).Value = 0.7103427268796919
$ws.Cells.Item(3, 17).Value = 3.687715201695778
$ws.Cells.Item(3, 18).Value = 33.189436815262
$ws.Cells.Item(3, 19).Value = 0.01855338466446163
$ws.Cells.Item(3, 20).Value = 0.01862985261096715

$ws.Cells.Item(4, 9).Value = 0.3422104954945279
$ws.Cells.Item(4, 10).Value = 0.3436209192170106
$ws.Cells.Item(4, 11).Value = 1
$ws.Cells.Item(4, 12).Value = 0.3333333333333333
$ws.Cells.Item(4, 13).Value = 0.1375686666666667
$ws.Cells.Item(4, 14).Value = 0.412706
$ws.Cells.Item(4, 15).Value = 0.2896572731203081
$ws.Cells.Item(4, 16).Value = 0.2896572731203081
$ws.Cells.Item(4, 17).Value = 19.70207589449556
$ws.Cells.Item(4, 18).Value = 177.31868305046
$ws.Cells.Item(4, 19).Value = 0.09912375895809446
$ws.Cells.Item(4, 20).Value = 0.09953229844749298

$ws.Cells.Item(5, 9).Value = 0.3422104954945279
$ws.Cells.Item(5, 10).Value = 0.3436209192170106
$ws.Cells.Item(5, 11).Value = 2
$ws.Cells.Item(5, 12).Value = 0.6666666666666666
$ws.Cells.Item(5, 13).Value = 0.3373673333333334
$ws.Cells.Item(5, 14).Value = 1.012102
$ws.Cells.Item(5, 15).Value = 0.7103427268796919
$ws.Cells.Item(5, 16).Value = 0.7103427268796919
$ws.Cells.Item(5, 17).Value = 48.31650234542445
$ws.Cells.Item(5, 18).Value = 434.84852110882
$ws.Cells.Item(5, 19).Value = 0.2430867365364335
$ws.Cells.Item(5, 20).Value = 0.2440886207695177

$ws.Cells.Item(6, 7).Value = 157.1889546666667
$ws.Cells.Item(6, 8).Value = 471.566864
$ws.Cells.Item(6, 9).Value = 0.3755976715691904
$ws.Cells.Item(6, 10).Value = 0.3771457008466821
$ws.Cells.Item(6, 11).Value = 1
$ws.Cells.Item(6, 12).Value = 0.3333333333333333
$ws.Cells.Item(6, 13).Value = 0.1375686666666667
$ws.Cells.Item(6, 14).Value = 0.412706
$ws.Cells.Item(6, 15).Value = 0.2896572731203081
$ws.Cells.Item(6, 16).Value = 0.2896572731203081
$ws.Cells.Item(6, 17).Value = 21.62427490822045
$ws.Cells.Item(6, 18).Value = 194.618474173984
$ws.Cells.Item(6, 19).Value = 0.1087945973370688
$ws.Cells.Item(6, 20).Value = 0.1092429952762974

$ws.Cells.Item(7, 7).Value = 157.1889546666667
$ws.Cells.Item(7, 8).Value = 471.566864
$ws.Cells.Item(7, 9).Value = 0.3755976715691904
$ws.Cells.Item(7, 10).Value = 0.3771457008466821
$ws.Cells.Item(7, 11).Value = 2
$ws.Cells.Item(7, 12).Value = 0.6666666666666666
$ws.Cells.Item(7, 13).Value = 0.3373673333333334
$ws.Cells.Item(7, 14).Value = 1.012102
$ws.Cells.Item(7, 15).Value = 0.7103427268796919
$ws.Cells.Item(7, 16).Value = 0.7103427268796919
$ws.Cells.Item(7, 17).Value = 53.03041846534756
$ws.Cells.Item(7, 18).Value = 477.273766188128
$ws.Cells.Item(7, 19).Value = 0.2668030742321216
$ws.Cells.Item(7, 20).Value = 0.2679027055703847

$ws.Cells.Item(8, 7).Value = 5.153359
$ws.Cells.Item(8, 8).Value = 10.306718
$ws.Cells.Item(8, 9).Value = 0.01231377640537609
$ws.Cells.Item(8, 10).Value = 0.008243018499152039
$ws.Cells.Item(8, 11).Value = 1
$ws.Cells.Item(8, 12).Value = 0.3333333333333333
$ws.Cells.Item(8, 13).Value = 0.1375686666666667
$ws.Cells.Item(8, 14).Value = 0.412706
$ws.Cells.Item(8, 15).Value = 0.2896572731203081
$ws.Cells.Item(8, 16).Value = 0.2896572731203081
$ws.Cells.Item(8, 17).Value = 0.7089407264846667
$ws.Cells.Item(8, 18).Value = 4.253644358908001
$ws.Cells.Item(8, 19).Value = 0.003566774895394428
$ws.Cells.Item(8, 20).Value = 0.002387650260744634

$ws.Cells.Item(9, 7).Value = 5.153359
$ws.Cells.Item(9, 8).Value = 10.306718
$ws.Cells.Item(9, 9).Value = 0.01231377640537609
$ws.Cells.Item(9, 10).Value = 0.008243018499152039
$ws.Cells.Item(9, 11).Value = 2
$ws.Cells.Item(9, 12).Value = 0.6666666666666666
$ws.Cells.Item(9, 13).Value = 0.3373673333333334
$ws.Cells.Item(9, 14).Value = 1.012102
$ws.Cells.Item(9, 15).Value = 0.7103427268796919
$ws.Cells.Item(9, 16).Value = 0.7103427268796919
$ws.Cells.Item(9, 17).Value = 1.738574983539333
$ws.Cells.Item(9, 18).Value = 10.431449901236
$ws.Cells.Item(9, 19).Value = 0.008747001509981662
$ws.Cells.Item(9, 20).Value = 0.005855368238407405

$ws.Cells.Item(10, 7).Value = 102.0140613333333
$ws.Cells.Item(10, 8).Value = 306.042184
$ws.Cells.Item(10, 9).Value = 0.2437591368004809
$ws.Cells.Item(10, 10).Value = 0.2447637923374727
$ws.Cells.Item(10, 11).Value = 1
$ws.Cells.Item(10, 12).Value = 0.3333333333333333
$ws.Cells.Item(10, 13).Value = 0.1375686666666667
$ws.Cells.Item(10, 14).Value = 0.412706
$ws.Cells.Item(10, 15).Value = 0.2896572731203081
$ws.Cells.Item(10, 16).Value = 0.2896572731203081
$ws.Cells.Item(10, 17).Value = 14.03393839887822
$ws.Cells.Item(10, 18).Value = 126.305445589904
$ws.Cells.Item(10, 19).Value = 0.07060660686378745
$ws.Cells.Item(10, 20).Value = 0.07089761264705771

$ws.Cells.Item(11, 7).Value = 102.0140613333333
$ws.Cells.Item(11, 8).Value = 306.042184
$ws.Cells.Item(11, 9).Value = 0.2437591368004809
$ws.Cells.Item(11, 10).Value = 0.2447637923374727
$ws.Cells.Item(11, 11).Value = 2
$ws.Cells.Item(11, 12).Value = 0.6666666666666666
$ws.Cells.Item(11, 13).Value = 0.3373673333333334
$ws.Cells.Item(11, 14).Value = 1.012102
$ws.Cells.Item(11, 15).Value = 0.7103427268796919
$ws.Cells.Item(11, 16).Value = 0.7103427268796919
$ws.Cells.Item(11, 17).Value = 34.41621183452978
$ws.Cells.Item(11, 18).Value = 309.7459065107681
$ws.Cells.Item(11, 19).Value = 0.1731525299366934
$ws.Cells.Item(11, 20).Value = 0.173866179690415
